$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoofdgegevens")

$ws.Range("B13").Value = 600
$ws.Range("B14").Formula = "=SUM(B15:B17)"
$ws.Range("B16").Value = 1000
$ws.Range("B17").Value = 10000000

$ws.Activate()
$ws.Range("F14").Select()
